$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell while forcing text storage
# (avoids "005571" -> 5571 / "0.91" -> 0.91-as-number coercion), then
# strip the number-format side effect the coercion-avoidance adds so the
# cell is left with no explicit style (matches cells that carry no `s`
# in the target sheet).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Step 1: rename the existing "总计" sheet to "2022-Q1" (keeps its
#     sheetId/rId), then duplicate it (whole-sheet Copy preserves
#     sheetPr/pageMargins/sheetFormatPr) right after itself to become the
#     new "总计" sheet with a fresh sheetId/rId -- matching the
#     workbook.xml diff exactly. ---
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"

# =========================================================================
# Sheet "2022-Q1": new fund-holding detail table (was the old 总计 sheet;
# its 6 rows of quarterly-summary content get replaced with the fresh
# per-fund breakdown, columns B..H headed, 3 data rows).
# =========================================================================

# Extend the header styling (bold/centered/bordered, style index 2) from
# the existing B1:D1 cells out to the new E1:H1 header cells.
$q1.Range("D1").Copy($q1.Range("E1"))
$q1.Range("D1").Copy($q1.Range("F1"))
$q1.Range("D1").Copy($q1.Range("G1"))
$q1.Range("D1").Copy($q1.Range("H1"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "005571"
$q1.Range("C2").Value = "中银证券新能源灵活配置混合A"
Set-TextValue $q1.Range("D2") "0.91"
Set-TextValue $q1.Range("E2") "90.25"
Set-TextValue $q1.Range("F2") "4.63"
Set-TextValue $q1.Range("G2") "0.0421"
$q1.Range("H2").Value = 7

# Row 3
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "005572"
$q1.Range("C3").Value = "中银证券新能源灵活配置混合C"
Set-TextValue $q1.Range("D3") "0.28"
Set-TextValue $q1.Range("E3") "90.25"
Set-TextValue $q1.Range("F3") "4.63"
Set-TextValue $q1.Range("G3") "0.0130"
$q1.Range("H3").Value = 7

# Row 4 (replaces the old row 4 content entirely)
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "003981"
$q1.Range("C4").Value = "中银证券瑞益灵活配置混合C"
Set-TextValue $q1.Range("D4") "0.21"
Set-TextValue $q1.Range("E4") "89.21"
Set-TextValue $q1.Range("F4") "2.84"
Set-TextValue $q1.Range("G4") "0.0060"
$q1.Range("H4").Value = 10

# Drop the leftover old rows 5 & 6 (former 2021-Q1 / 2020-Q4 summary rows).
$q1.Range("A5:A6").EntireRow.Delete()

# =========================================================================
# Sheet "总计": quarterly summary table, now with a new 2022-Q1 row
# inserted at the top (everything else shifts down by one row).
# =========================================================================

# Insert a fresh row 2 (pushes the old rows 2-6 down to 3-7) and strip the
# formatting it inherits from the header row above so it starts blank.
$total.Range("A2").EntireRow.Insert()
$total.Range("B2:D2").ClearFormats()

# Give the new A2 the same bold/centered/bordered style as the other
# index cells in column A, then fill in the new top row's data.
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.06

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# --- Restore the original active sheet (2020-Q4) so `tabSelected` ends
#     up back where it started rather than on the newly-added sheet. ---
$wb.Worksheets.Item("2020-Q4").Activate()
